$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for numeric-looking price (D) and hour (G) values
# by forcing Text number format before assigning, so Excel does not
# silently convert these strings into numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "244.20"
$ws.Range("G2").Value = "19"

$ws.Range("D3").Value = "23.98"
$ws.Range("G3").Value = "19"

$ws.Range("D4").Value = "5.255"
$ws.Range("G4").Value = "19"

$ws.Range("D5").Value = "0.05828"
$ws.Range("G5").Value = "19"

$ws.Range("G6").Value = "19"

$ws.Range("D7").Value = "3.233"
$ws.Range("G7").Value = "19"

$ws.Range("D8").Value = "0.8080"
$ws.Range("G8").Value = "19"

$ws.Range("D9").Value = "0.8842"
$ws.Range("G9").Value = "19"

$ws.Range("D10").Value = "0.1383"
$ws.Range("G10").Value = "19"

$ws.Range("D11").Value = "0.07138"
$ws.Range("G11").Value = "19"

$ws.Range("D12").Value = "0.03080"
$ws.Range("G12").Value = "19"

$ws.Range("D13").Value = "0.03033"
$ws.Range("G13").Value = "19"

$ws.Range("D14").Value = "0.09335"
$ws.Range("G14").Value = "19"

$ws.Range("D15").Value = "3.821"
$ws.Range("G15").Value = "19"

$ws.Range("D16").Value = "0.001554"
$ws.Range("G16").Value = "19"

$ws.Range("D17").Value = "0.04709"
$ws.Range("G17").Value = "19"

$ws.Range("D18").Value = "0.0006049"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").Value = "19"

$ws.Range("D19").Value = "0.006212"
$ws.Range("G19").Value = "19"

$ws.Range("D20").Value = "0.001260"
$ws.Range("G20").Value = "19"

$ws.Range("D21").Value = "0.004078"
$ws.Range("G21").Value = "19"

$ws.Range("D22").Value = "0.00008704"
$ws.Range("G22").Value = "19"

$ws.Range("G23").Value = "19"

$ws.Range("D24").Value = "2.163"
$ws.Range("G24").Value = "19"

$ws.Range("D25").Value = "0.3183"
$ws.Range("G25").Value = "19"

$ws.Range("D26").Value = "0.1314"
$ws.Range("G26").Value = "19"

$ws.Range("G27").Value = "19"

$ws.Range("D28").Value = "0.0001970"
$ws.Range("G28").Value = "19"

$ws.Range("G29").Value = "19"

$ws.Range("G30").Value = "19"

$ws.Range("G31").Value = "19"

$ws.Range("G32").Value = "19"

$ws.Range("G33").Value = "19"

$ws.Range("G34").Value = "19"

$ws.Range("G35").Value = "19"

$ws.Range("G36").Value = "19"

$ws.Range("G37").Value = "19"

$ws.Range("G38").Value = "19"

$ws.Range("G39").Value = "19"

$ws.Range("D40").Value = "0.03847"
$ws.Range("G40").Value = "19"

$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1055"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("G41").Value = "19"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002542"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("G42").Value = "19"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.003238"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"
$ws.Range("G43").Value = "19"

$ws.Range("D44").Value = "0.007246"
$ws.Range("G44").Value = "19"

$ws.Range("D45").Value = "0.00005343"
$ws.Range("G45").Value = "19"

$ws.Range("G46").Value = "19"

$ws.Range("D47").Value = "0.4899"
$ws.Range("G47").Value = "19"

$ws.Range("D48").Value = "0.002240"
$ws.Range("G48").Value = "19"

$ws.Range("G49").Value = "19"

$ws.Range("G50").Value = "19"

$ws.Range("G51").Value = "19"
